$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking price strings to stay as text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "64.358.41"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.138.85"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "571.70"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "164.08"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -5.64%  "
$ws.Range("D9").Value = "3.152.28"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").Value = "6.62"
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").Value = "0.384"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "3.687.25"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "64.353.43"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "24.98"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "3.141.63"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "405.79"
$ws.Range("E19").Value = "  -3.01%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "68.96"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E27").Value = "  -3.72%  "
$ws.Range("D28").Value = "8.88"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").Value = "162.79"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "1.36"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "2.634.34"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("D40").Value = "23.65"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "38.29"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").Value = "0.691"
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("D44").Value = "0.0610"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").Value = "5.38"
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("D46").Value = "0.0255"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("D47").Value = "287.75"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "21.20"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("E51").Value = "  +0.20%  "
